# Apply the OOXML diff:
#  1. Remove the custom paragraph style "AbstractTitle" entirely.
#  2. Change the "Abstract" style's paragraph spacing "before" value
#     from 100 twips (5pt) to 300 twips (15pt).

$d = $word.ActiveDocument

# 1) Delete the "Abstract Title" style (styleId="AbstractTitle").
$abstractTitle = $d.Styles.Item("Abstract Title")
$abstractTitle.Delete()

# 2) Update the "Abstract" style's space-before from 5pt (100 twips) to
#    15pt (300 twips); SpaceBefore/SpaceAfter are expressed in points.
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15
